# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates the td_sim_1 (column C) simulation values for rows 2-97 and
# recomputes the record_atd (column D) midpoint values plus the
# average_simulation_TD summary cell (C98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New td_sim_1 (column C) values for rows 2 through 97.
$newC = @(
    411,25,322,76,39,243,136,158,17,94,
    115,68,131,72,211,75,36,35,345,9,
    40,21,50,99,41,31,71,13,38,186,
    332,146,44,90,24,1097,486,451,223,325,
    240,6,8,49,341,15,274,249,142,54,
    434,12,112,19,334,28,308,111,65,20,
    471,276,119,66,95,167,7,42,58,47,
    190,61,113,138,445,23,100,80,110,444,
    77,437,443,11,453,164,267,73,51,482,
    48,212,14,571,123,69
)

$firstRow = 2
$lastRow = 97

for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $firstRow + $i
    $cCell = $ws.Cells.Item($row, 3)   # column C
    $dCell = $ws.Cells.Item($row, 4)   # column D
    $bCell = $ws.Cells.Item($row, 2)   # column B

    $cCell.Value2 = $newC[$i]

    $bVal = $bCell.Value2
    if ($null -eq $bVal) {
        # No td_sim_0 value recorded for this record -> record_atd equals td_sim_1
        $dCell.Value2 = $newC[$i]
    } else {
        $dCell.Value2 = ($bVal + $newC[$i]) / 2
    }
}

# Recompute the average_simulation_TD summary value in C98 as the mean
# of the updated td_sim_1 column (C2:C97).
$sumC = 0.0
for ($i = 0; $i -lt $newC.Length; $i++) {
    $sumC += $newC[$i]
}
$ws.Cells.Item(98, 3).Value2 = $sumC / $newC.Length
